$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.02258322285507441
$ws.Range("C2").Value = 0.05231270169004087
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 250.1427460453503

$ws.Range("B3").Value = 0.02258322285507441
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 1.572787616952587

$ws.Range("B4").Value = 3.182878228561681
$ws.Range("C4").Value = 86.29678392075563
$ws.Range("D4").Value = 3.082599426703578
$ws.Range("E4").Value = 246.9852506941017
$ws.Range("G4").Value = 339.5475122701225

$ws.Range("B5").Value = 0.02258322285507441
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 3.082599426703578
$ws.Range("E5").Value = 246.9852506941017
$ws.Range("G5").Value = 251.7436698025591

$ws.Range("B6").Value = 0.06328177979961902
$ws.Range("C6").Value = 9.226618575922256
$ws.Range("D6").Value = 3.082599426703578
$ws.Range("E6").Value = 6.48142807727062
$ws.Range("G6").Value = 18.85392785969607

$ws.Range("B7").Value = 0.3464964993005633
$ws.Range("C7").Value = 1.65323645889881
$ws.Range("D7").Value = 0.1529057820181812
$ws.Range("E7").Value = 246.9852506941017
$ws.Range("G7").Value = 249.1378894343192

$ws.Range("B8").Value = 1.505614041169197
$ws.Range("C8").Value = 1.65323645889881
$ws.Range("D8").Value = 3.082599426703578
$ws.Range("E8").Value = 0.4998867070740569
$ws.Range("G8").Value = 6.741336633845642

$ws.Range("B9").Value = 0.1554434735375247
$ws.Range("C9").Value = 0.3375848360084654
$ws.Range("D9").Value = 0.7127328510149897
$ws.Range("E9").Value = 6.48142807727062
$ws.Range("G9").Value = 7.6871892378316

$ws.Range("B10").Value = 1.505614041169197
$ws.Range("C10").Value = 1.65323645889881
$ws.Range("D10").Value = 0.7127328510149897
$ws.Range("E10").Value = 0.4998867070740569
$ws.Range("G10").Value = 4.371470058157054

$ws.Range("B11").Value = 3.182878228561681
$ws.Range("C11").Value = 1.65323645889881
$ws.Range("D11").Value = 0.7127328510149897
$ws.Range("E11").Value = 71517.89157740913
$ws.Range("G11").Value = 71523.44042494761

$ws.Range("B12").Value = 0.1554434735375247
$ws.Range("C12").Value = 0.004309184025731883
$ws.Range("D12").Value = 16.98373111632243
$ws.Range("E12").Value = 71517.89157740913
$ws.Range("G12").Value = 71535.03506118302

$ws.Range("B13").Value = 3.182878228561681
$ws.Range("C13").Value = 1.65323645889881
$ws.Range("D13").Value = 0.7127328510149897
$ws.Range("E13").Value = 0.4998867070740569
$ws.Range("G13").Value = 6.048734245549538

$ws.Range("B14").Value = 3.182878228561681
$ws.Range("C14").Value = 1.65323645889881
$ws.Range("D14").Value = 16.98373111632243
$ws.Range("E14").Value = 6.48142807727062
$ws.Range("G14").Value = 28.30127388105354

